$wb = $excel.ActiveWorkbook

# Avoid the "are you sure you want to delete" prompt when removing a sheet.
$excel.DisplayAlerts = $false

# Rename "Paineis DARQ" -> "PAINEIS DARQ"
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"

# Rename "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# Remove the "Desarquivamentos Pendentes" sheet entirely
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()

$excel.DisplayAlerts = $true
